$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26 (the Total row), shifting it down to row 27
$ws.Rows.Item(26).Insert()

# Copy formatting from the single cells in row 24 (the last data row) into row 25
$ws.Cells.Item(24, 1).Copy()
$ws.Cells.Item(25, 1).PasteSpecial(-4122)
$ws.Cells.Item(24, 2).Copy()
$ws.Cells.Item(25, 2).PasteSpecial(-4122)
$ws.Cells.Item(24, 3).Copy()
$ws.Cells.Item(25, 3).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row 25 data
$ws.Cells.Item(25, 1).Value = "Completed Home and About pages"
$ws.Cells.Item(25, 2).Value = 43374
$ws.Cells.Item(25, 3).Value = 0.5

# Fix up the Total row formula to include the new row
$ws.Cells.Item(27, 3).Formula = "=SUM(C2:C26)"

# Update the selection to match
$ws.Range("C26").Select()
